$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, `
                             $true, 1, $false, $new, 2) | Out-Null
}

Replace-Text "2025-08-20 Wednesday" "2025-08-21 Thursday"

Replace-Text "66×40=" "55×89="
Replace-Text "98×29=" "84×25="
Replace-Text "51×70=" "26×91="
Replace-Text "82×82=" "83×39="
Replace-Text "41×59=" "43×56="

Replace-Text "85×73=" "53×62="
Replace-Text "94×34=" "14×47="
Replace-Text "38×55=" "79×63="
Replace-Text "46×36=" "31×42="
Replace-Text "61×66=" "55×21="

Replace-Text "73×36=" "71×13="
Replace-Text "94×42=" "36×72="
Replace-Text "98×66=" "80×78="
Replace-Text "67×50=" "38×68="
Replace-Text "27×65=" "54×32="

Replace-Text "83×78=" "80×88="
Replace-Text "74×34=" "88×58="
Replace-Text "86×91=" "54×68="
Replace-Text "23×39=" "11×90="
Replace-Text "87×66=" "59×56="

Replace-Text "43×17=" "87×83="
Replace-Text "13×79=" "39×74="
Replace-Text "45×21=" "27×31="
Replace-Text "97×69=" "82×41="
Replace-Text "16×91=" "61×68="
